# "1-100 compare number of quanv filter on mnist fashion"
#
# Sheet "Diff num of quanv filter": drop the standalone "MNIST" title row,
# add a second comparison column C ("mnist_fashion") next to the existing
# "mnist" column B, fill in the updated B-column values, and relocate the
# two little color/author legend blocks that used to sit in columns D/E up
# against the new layout (D:E14-16 -> F:G9-11, D18-22 -> F13-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diff num of quanv filter")

# --- Drop the old "MNIST" title row and the leftover "Status" label ---
$ws.Range("A2").ClearContents()
$ws.Range("D6").ClearContents()

# --- Header row: title stays, new mnist / mnist_fashion column headers ---
$ws.Range("B4").Value = "mnist"
$ws.Range("C4").Value = "mnist_fashion"

# --- New column C ("mnist_fashion"), styled like column B (yellow, s=3) ---
$cValues = @{5=2;6=2;7=2;8=2;9=2;10=1;11=2;12=2;13=2;14=1;15=1;16=1;17=1;18=1;19=1}
foreach ($r in 5..19) {
  $dst = $ws.Cells.Item($r, 3)
  $ws.Range("B6").Copy($dst)
  $dst.Value = $cValues[$r]
}

# --- Existing column B ("mnist") values that changed to "2" ---
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 2

# B12/B13 previously had no value (only the "no fill" style s=5); B14 was
# still "no fill" too. All three now take the yellow fill (s=3), like B7.
foreach ($r in 12..14) {
  $dst = $ws.Cells.Item($r, 2)
  $ws.Range("B7").Copy($dst)
  $dst.Value = 2
}

# --- Move the status-code legend block (D14:E16 -> F9:G11) ---
$pairs = @(
  @("D14:E14", "F9"),
  @("D15:E15", "F10"),
  @("D16:E16", "F11")
)
foreach ($p in $pairs) {
  $src = $ws.Range($p[0])
  $src.Cut($ws.Range($p[1]))
  $src.Clear()
}

# --- Move the name/location legend block (D18:D22 -> F13:F17) ---
foreach ($r in 18..22) {
  $src = $ws.Cells.Item($r, 4)
  $dst = $ws.Cells.Item($r - 5, 6)
  $src.Cut($dst)
  $src.Clear()
}

# --- Sheet view: selection now sits on the newly-added C19 cell ---
$ws.Range("C19").Select()
